$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B:E keep their original text storage (avoid COM auto-converting
# numeric-looking strings like "1.00" into the number 1, or URLs/names into other types).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '63.161.43'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '3.392.96'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '566.49'
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").Value = '155.21'
$ws.Range("E6").Value = '  +1.80%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '3.391.87'
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("D9").Value = '0.543'
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("E10").Value = '  -0.73%  '
$ws.Range("E11").Value = '  +2.69%  '
$ws.Range("D12").Value = '0.431'
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = '3.980.71'
$ws.Range("E13").Value = '  +1.24%  '
$ws.Range("E14").Value = '  -3.18%  '
$ws.Range("D15").Value = '0.0000191'
$ws.Range("E15").Value = '  +6.33%  '
$ws.Range("D16").Value = '27.15'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").Value = '63.218.13'
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D18").Value = '3.355.64'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").Value = '6.24'
$ws.Range("E19").Value = '  -2.02%  '
$ws.Range("E20").Value = '  +1.57%  '
$ws.Range("D21").Value = '377.93'
$ws.Range("E22").Value = '  -3.55%  '
$ws.Range("D23").Value = '0.996'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").Value = '71.38'
$ws.Range("E24").Value = '  +1.43%  '
$ws.Range("E25").Value = '  -1.72%  '
$ws.Range("E26").Value = '  +24.38%  '
$ws.Range("D27").Value = '9.41'
$ws.Range("E27").Value = '  +6.08%  '
$ws.Range("D28").Value = '0.178'
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = '6.05'
$ws.Range("E30").Value = '  +7.53%  '
$ws.Range("D31").Value = '1.35'
$ws.Range("E31").Value = '  +3.70%  '
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").Value = '23.13'
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").Value = '6.33'
$ws.Range("E35").Value = '  -3.88%  '
$ws.Range("D36").Value = '6.78'
$ws.Range("E36").Value = '  +1.55%  '
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("D39").Value = '2.953.87'
$ws.Range("E39").Value = '  +5.00%  '
$ws.Range("D40").Value = '27.02'
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("D41").Value = '0.0752'
$ws.Range("E41").Value = '  +1.17%  '
$ws.Range("D43").Value = '0.0317'
$ws.Range("E43").Value = '  +2.01%  '
$ws.Range("D44").Value = '41.61'
$ws.Range("E44").Value = '  +2.38%  '
$ws.Range("D45").Value = '0.760'
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").Value = '23.20'
$ws.Range("E47").Value = '  +5.82%  '
$ws.Range("E48").Value = '  +3.32%  '
$ws.Range("E49").Value = '  +20.53%  '
$ws.Range("D50").Value = '6.34'
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("D51").Value = '0.829'
$ws.Range("E51").Value = '  +3.77%  '
